$paraData = @(
    @{ bold = $true; color = $null; segs = @() },
    @{ bold = $false; color = $null; segs = @(@{ t = "text"; v = "SAT Dec 02" }, @{ t = "text"; v = " 10:20:57 PST 2017" }) },
    @{ bold = $false; color = $null; segs = @(@{ t = "text"; v = "Person Name" }, @{ t = "tab"; v = $null }, @{ t = "tab"; v = $null }, @{ t = "tab"; v = $null }, @{ t = "tab"; v = $null }, @{ t = "text"; v = "- BASAVA" }) },
    @{ bold = $false; color = $null; segs = @(@{ t = "text"; v = "Bill number" }, @{ t = "tab"; v = $null }, @{ t = "tab"; v = $null }, @{ t = "tab"; v = $null }, @{ t = "tab"; v = $null }, @{ t = "text"; v = "- 1735" }) },
    @{ bold = $false; color = $null; segs = @(@{ t = "text"; v = "---------------------------------------------------------------" }) },
    @{ bold = $false; color = $null; segs = @(@{ t = "text"; v = "Item Name" }, @{ t = "tab"; v = $null }, @{ t = "tab"; v = $null }, @{ t = "tab"; v = $null }, @{ t = "tab"; v = $null }, @{ t = "text"; v = "- POTATO" }) },
    @{ bold = $false; color = $null; segs = @(@{ t = "text"; v = "Number of Pockets" }, @{ t = "tab"; v = $null }, @{ t = "tab"; v = $null }, @{ t = "tab"; v = $null }, @{ t = "text"; v = "- 4" }) },
    @{ bold = $false; color = $null; segs = @(@{ t = "text"; v = "Number of KGs" }, @{ t = "tab"; v = $null }, @{ t = "tab"; v = $null }, @{ t = "tab"; v = $null }, @{ t = "text"; v = "- 198" }) },
    @{ bold = $false; color = $null; segs = @(@{ t = "text"; v = "Rate" }, @{ t = "tab"; v = $null }, @{ t = "tab"; v = $null }, @{ t = "tab"; v = $null }, @{ t = "tab"; v = $null }, @{ t = "tab"; v = $null }, @{ t = "text"; v = "- 12" }) },
    @{ bold = $false; color = $null; segs = @(@{ t = "text"; v = "Total Price" }, @{ t = "tab"; v = $null }, @{ t = "tab"; v = $null }, @{ t = "tab"; v = $null }, @{ t = "tab"; v = $null }, @{ t = "text"; v = "- 2376.0" }) },
    @{ bold = $false; color = "FF0000"; segs = @(@{ t = "text"; v = "Amount Received" }, @{ t = "tab"; v = $null }, @{ t = "tab"; v = $null }, @{ t = "tab"; v = $null }, @{ t = "text"; v = "- 2000" }) },
    @{ bold = $false; color = $null; segs = @(@{ t = "text"; v = "Amount balance" }, @{ t = "tab"; v = $null }, @{ t = "tab"; v = $null }, @{ t = "tab"; v = $null }, @{ t = "text"; v = "- 28323.0" }) },
    @{ bold = $false; color = $null; segs = @(@{ t = "text"; v = "Amount Received mode" }, @{ t = "tab"; v = $null }, @{ t = "tab"; v = $null }, @{ t = "text"; v = "- CASH" }) },
    @{ bold = $false; color = $null; segs = @() },
    @{ bold = $false; color = $null; segs = @(@{ t = "text"; v = "Item Name" }, @{ t = "tab"; v = $null }, @{ t = "tab"; v = $null }, @{ t = "tab"; v = $null }, @{ t = "tab"; v = $null }, @{ t = "text"; v = "- CARROT" }) },
    @{ bold = $false; color = $null; segs = @(@{ t = "text"; v = "Number of Pockets" }, @{ t = "tab"; v = $null }, @{ t = "tab"; v = $null }, @{ t = "tab"; v = $null }, @{ t = "text"; v = "- 1" }) },
    @{ bold = $false; color = $null; segs = @(@{ t = "text"; v = "Number of KGs" }, @{ t = "tab"; v = $null }, @{ t = "tab"; v = $null }, @{ t = "tab"; v = $null }, @{ t = "text"; v = "- 94" }) },
    @{ bold = $false; color = $null; segs = @(@{ t = "text"; v = "Rate" }, @{ t = "tab"; v = $null }, @{ t = "tab"; v = $null }, @{ t = "tab"; v = $null }, @{ t = "tab"; v = $null }, @{ t = "tab"; v = $null }, @{ t = "text"; v = "- 58" }) },
    @{ bold = $false; color = $null; segs = @(@{ t = "text"; v = "Total Price" }, @{ t = "tab"; v = $null }, @{ t = "tab"; v = $null }, @{ t = "tab"; v = $null }, @{ t = "tab"; v = $null }, @{ t = "text"; v = "- 5452.0" }) },
    @{ bold = $true; color = $null; segs = @(@{ t = "text"; v = "Amount balance" }, @{ t = "tab"; v = $null }, @{ t = "tab"; v = $null }, @{ t = "tab"; v = $null }, @{ t = "text"; v = "- 33775.0" }) },
    @{ bold = $false; color = $null; segs = @() },
    @{ bold = $true; color = $null; segs = @() },
)
$d = $word.ActiveDocument

# ---- Step 1: merge the "THU NOV 30" / " 10:03:58 PST 2017" runs into one run ----
$null = $d.Content.Find.Execute("THU NOV 30 10:03:58 PST 2017", $false, $false, $false, $false, $false, $true, 1, $false, "THU NOV 30 10:03:58 PST 2017", 2)

# ---- Step 2: locate anchor paragraph ("Amount balance ... - 27947.0") ----
$anchorIdx = 0
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $t = $d.Paragraphs($i).Range.Text
    if ($t -like "*27947.0*") {
        $anchorIdx = $i
        break
    }
}

$anchorPara = $d.Paragraphs($anchorIdx)
$cursor = $anchorPara.Range.End - 1

$placeholders = @()

foreach ($p in $paraData) {
    # Insert paragraph break
    $ins = $d.Range($cursor, $cursor)
    $ins.InsertAfter("`r")
    $cursor = $ins.End

    if ($p.segs.Count -eq 0) {
        # Empty paragraph - insert placeholder char, format it, note location for later deletion
        $ins2 = $d.Range($cursor, $cursor)
        $ins2.InsertAfter("@")
        $phStart = $cursor
        $phEnd = $ins2.End
        $cursor = $phEnd
        $phRange = $d.Range($phStart, $phEnd)
        if ($p.bold) {
            $phRange.Font.Bold = 1
        }
        if ($p.color) {
            $phRange.Font.Color = [long]("0x" + $p.color.Substring(4,2) + $p.color.Substring(2,2) + $p.color.Substring(0,2))
        }
        $placeholders += $phStart
    } else {
        $paraTextStart = $cursor
        $prevType = $null
        for ($si = 0; $si -lt $p.segs.Count; $si++) {
            $seg = $p.segs[$si]
            $ins3 = $d.Range($cursor, $cursor)
            if ($seg.t -eq "tab") {
                $ins3.InsertAfter("`t")
            } else {
                $ins3.InsertAfter($seg.v)
            }
            $segEnd = $ins3.End
            # Force a separate run when two consecutive "text" segments would
            # otherwise merge (identical formatting) - toggle bold on/off.
            if ($seg.t -eq "text" -and $prevType -eq "text") {
                $segRange = $d.Range($cursor, $segEnd)
                $segRange.Font.Bold = 1
                $segRange.Font.Bold = 0
            }
            $cursor = $segEnd
            $prevType = $seg.t
        }
        $paraRange = $d.Range($paraTextStart, $cursor)
        if ($p.bold) {
            $paraRange.Font.Bold = 1
        }
        if ($p.color) {
            $paraRange.Font.Color = [long]("0x" + $p.color.Substring(4,2) + $p.color.Substring(2,2) + $p.color.Substring(0,2))
        }
    }
}

Write-Output ("Done inserting. Paragraphs.Count=" + $d.Paragraphs.Count)
Write-Output ("Placeholders=" + ($placeholders -join ","))

# ---- Step 3: remove placeholder chars, in reverse order so earlier offsets stay valid ----
for ($i = $placeholders.Count - 1; $i -ge 0; $i--) {
    $pos = $placeholders[$i]
    $delRange = $d.Range($pos, $pos + 1)
    $delRange.Delete()
}

Write-Output ("Final Paragraphs.Count=" + $d.Paragraphs.Count)
